$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.240107297897339
$ws.Range("B1").Value = 2.890908241271973
$ws.Range("C1").Value = 2.542561054229736
$ws.Range("D1").Value = 2.765085458755493
$ws.Range("E1").Value = 3.055150747299194
